$p = $ppt.ActivePresentation

# Slide 2 ("// summary: presentation") - capitalize agenda bullet points
# and rename the last bullet from "demo of Pokemon page" to "Demo of Game".
$s2 = $p.Slides.Item(2)
$contentShape = $s2.Shapes.Item(2)
$tr = $contentShape.TextFrame.TextRange

$tr.Paragraphs(1, 1).Runs(1, 1).Text = "Summary of task"
$tr.Paragraphs(2, 1).Runs(1, 1).Text = "Summary of the week"
$tr.Paragraphs(3, 1).Runs(1, 1).Text = "Code snippets"
$tr.Paragraphs(4, 1).Runs(1, 1).Text = "Lessons learned"
$tr.Paragraphs(5, 1).Runs(1, 1).Text = "Demo of Game"

# Slide 5 ("// snippets") - nudge the horizontal divider line's width by
# 1 EMU (10515598 -> 10515597), matching a resave rounding drift.
$s5 = $p.Slides.Item(5)
$line = $s5.Shapes.Item(5)
$line.Width = 827.999842519685
